# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# Diff summary: cell B11 on the "Rules" sheet changes from the shared
# string "R40" to the (new) shared string "1". The cell's existing
# style/formatting is left untouched - this mirrors a user typing the
# digit 1 into a General-formatted cell but forcing it to be stored as
# text (the leading apostrophe is Excel's standard "treat as text"
# quote-prefix convention), rather than letting it parse as the number 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Cells.Item(11, 2).Value = "'1"
